# CIERRE 26 DIC 23
# Record two new water-bottle purchase entries (rows 71 and 72) that were
# previously blank placeholder rows, and move the active selection down to
# the next empty row (D73), matching the ledger's running-balance pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 71: 22-Dec-2023, "3 botellones", -159
$ws.Range("B71").Value = 45282
$ws.Range("C71").Value = "3 botellones"
$ws.Range("D71").Value = -159

# Row 72: 26-Dec-2023, "4 botellones", -212
$ws.Range("B72").Value = 45286
$ws.Range("C72").Value = "4 botellones"
$ws.Range("D72").Value = -212

# Move selection to the next blank entry row, as in the closed-out ledger.
$ws.Range("D73").Select()
